# Fruta / hortaliza, semanal
# The source rows (2-5) get cyclically rotated: data that was in row 5
# moves to row 2, row 4's data moves to row 3, row 2's data moves to
# row 4, and row 3's data moves to row 5. Only columns D and J..Q carry
# data that differs row to row (A,B,C,E,F,G,H,I,R are identical across
# all rows), so we only need to touch those.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "before" values for the columns that vary across rows.
$rows = 2..5
$data = @{}
foreach ($r in $rows) {
    $data[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        J = $ws.Cells.Item($r, 10).Value2
        K = $ws.Cells.Item($r, 11).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
        Q = $ws.Cells.Item($r, 17).Value2
    }
}

# new_row <- old_row mapping observed in the target workbook.
$mapping = @{ 2 = 5; 3 = 4; 4 = 2; 5 = 3 }

foreach ($newRow in $rows) {
    $oldRow = $mapping[$newRow]
    $src = $data[$oldRow]

    $ws.Cells.Item($newRow, 4).Value2 = $src.D
    $ws.Cells.Item($newRow, 10).Value2 = $src.J
    $ws.Cells.Item($newRow, 11).Value2 = $src.K
    $ws.Cells.Item($newRow, 12).Value2 = $src.L
    $ws.Cells.Item($newRow, 13).Value2 = $src.M
    $ws.Cells.Item($newRow, 14).Value2 = $src.N
    $ws.Cells.Item($newRow, 15).Value2 = $src.O
    $ws.Cells.Item($newRow, 16).Value2 = $src.P
    $ws.Cells.Item($newRow, 17).Value2 = $src.Q
}
